$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: "Edges" (was E1) moves into C1 (was "Yc"); D1 ("Zc") and E1 are removed ---
$ws.Range("C1").Value = "Edges"
$ws.Range("D1").ClearContents()
$ws.Range("E1").ClearContents()

# --- Row 2: new data values ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 0.125
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 5

# --- Row 3: clear all previous data, leave only a formatted (quote-prefix) empty cell at F3 ---
$ws.Range("H3").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A3:E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("H3").Clear()

# --- Column widths: drop the narrow column C (6.71) and set column C to 15.14 (was column E) ---
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# --- Selection ---
$ws.Range("H7").Select() | Out-Null
